$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Luca"
$ws.Range("B14").Value = "GDPR"
$ws.Range("C14").Value = 105

$ws.Range("D13").Copy()
$ws.Range("D14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D14").Value = 43523

$ws.Range("A15").Select()
